$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("M:O").Insert()
$ws.Columns("R:T").Insert()

$ws.Range("H1").Value2 = "arrecadado_avg"
$ws.Range("I1").Value2 = "arrecadado_std"
$ws.Range("J1").Value2 = "arrecadado_min"
$ws.Range("K1").Value2 = "arrecadado_max"

$ws.Range("M1").Value2 = "apoio_std"
$ws.Range("N1").Value2 = "apoio_min"
$ws.Range("O1").Value2 = "apoio_max"

$ws.Range("Q1").Value2 = "contribuicoes_med"
$ws.Range("R1").Value2 = "contribuicoes_std"
$ws.Range("S1").Value2 = "contribuicoes_min"
$ws.Range("T1").Value2 = "contribuicoes_max"

$ws.Range("M2:O6").NumberFormat = "R$ #,##0.00"
$ws.Range("R2:T6").NumberFormat = "#,##0"

$ws.Range("L2").Value2 = 31.16847126718795
$ws.Range("M2").Value2 = 26.9469146898807
$ws.Range("N2").Value2 = 8.140546434454963
$ws.Range("O2").Value2 = 84.0771316599004

$ws.Range("L3").Value2 = 42.14013096402113
$ws.Range("M3").Value2 = 8.830628986869351
$ws.Range("N3").Value2 = 35.89593332526331
$ws.Range("O3").Value2 = 48.38432860277894

$ws.Range("L4").Value2 = 19.42257389357928
$ws.Range("M4").Value2 = 8.876706034650423
$ws.Range("N4").Value2 = 5.929916345397809
$ws.Range("O4").Value2 = 35.80030877323957

$ws.Range("L5").Value2 = 17.81312171425239
$ws.Range("M5").Value2 = 9.871079671113662
$ws.Range("N5").Value2 = 6.098311514417047
$ws.Range("O5").Value2 = 45.46067338136409

$ws.Range("L6").Value2 = 21.37695663886886
$ws.Range("M6").Value2 = 15.58070588764584
$ws.Range("N6").Value2 = 1.011042153300025
$ws.Range("O6").Value2 = 70.01644246718027

$ws.Range("R2").Value2 = 2.497617912751115
$ws.Range("S2").Value2 = 3
$ws.Range("T2").Value2 = 10

$ws.Range("R3").Value2 = 3.535533905932738
$ws.Range("S3").Value2 = 10
$ws.Range("T3").Value2 = 15

$ws.Range("R4").Value2 = 20.56688435388656
$ws.Range("S4").Value2 = 1
$ws.Range("T4").Value2 = 79

$ws.Range("R5").Value2 = 45.19114957599552
$ws.Range("S5").Value2 = 1
$ws.Range("T5").Value2 = 208

$ws.Range("R6").Value2 = 30.79468667274807
$ws.Range("S6").Value2 = 1
$ws.Range("T6").Value2 = 196
